{"js": "// Insert a new \"List Bullet\" paragraph with the new docente's name right\n// after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) =>\n  p.text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1\n);\n\nif (!target) {\n  throw new Error('Could not find \"Docente(s) Respons\u00e1vel(eis)\" paragraph');\n}\n\nconst newPara = target.insertParagraph(\"4893449 - D\u00e9bora Souza Alvim\", \"After\");\nnewPara.style = \"List Bullet\";\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph with the new docente's name right\n# after the \"Docente(s) Respons\u00e1vel(eis)\" heading paragraph.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like '*Docente(s) Respons\u00e1vel(eis)*') {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.Text = \"4893449 - D\u00e9bora Souza Alvim\"\n$newPara.Style = \"List Bullet\"\n"}
